$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray duplicated intro paragraph in "Problema 6"
#    ("Costruite un programma per gestire una lista di film. Per
#    immagazzinare i film dovete usare la seguente struttura" - the one
#    WITHOUT the trailing colon, right after the "Problema 6" heading).
#    The whole paragraph (text + its own paragraph mark) must go away.
# ---------------------------------------------------------------------
$targetText = "Costruite un programma per gestire una lista di film. Per immagazzinare i film dovete usare la seguente struttura" + [char]13
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq $targetText) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# ---------------------------------------------------------------------
# 2)-4) Collapse the "<w:tab/></w:r><w:r>struct </w:r><w:r>nodo..." runs
#    into a single run "<w:tab/><w:t>struct nodo..._t ...;</w:t>" for the
#    three struct-member declaration lines.
# ---------------------------------------------------------------------
function Merge-StructRun($fullLine) {
    # $fullLine example: "`tstruct nodo_t *next;"
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $full = $p.Range
        $text = $full.Text
        if ($text -eq ($fullLine + [char]13)) {
            $r = $d.Range($full.Start, $full.End - 1)
            $afterTab = $fullLine.Substring(1)
            $xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr/><w:tab/><w:t>$afterTab</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
            $r.InsertXML($xmlFrag)
            return $true
        }
    }
    return $false
}

Merge-StructRun ([char]9 + "struct nodo_t *next;") | Out-Null
Merge-StructRun ([char]9 + "struct nodo_t  *down;") | Out-Null
Merge-StructRun ([char]9 + "struct nodo_f_t *next;") | Out-Null

Write-Output "done"
